# Deploying to gh-pages from @ Alvearie/alvearie-fhir-ig@8e4a450c507ef6f746e072652acbb72e9504f19a
#
# Updates the StructureDefinition-rx-count.xlsx metadata & element tables:
#  - Sheet "Metadata": bump Version 5.0.0 -> 6.0.0, bump Date, fill in
#    Publisher, and replace the stray duplicated "Contact" row with a
#    "Jurisdiction" row (which removes one row from the table).
#  - Sheet "Elements": give the root Extension row a real Short/Definition
#    instead of the generic placeholder text.

$wb = $excel.ActiveWorkbook

$meta = $wb.Worksheets.Item(1)

# Version
$meta.Range("B3").Value = "6.0.0"

# Date
$meta.Range("B8").Value = "2022-01-21T20:46:54+00:00"

# Publisher
$meta.Range("B9").Value = "Alvearie Team"

# Replace the first "Contact" row with the new "Jurisdiction" row ...
$meta.Range("A10").Value = "Jurisdiction"
$meta.Range("B10").Value = "United States of America"

# ... and drop the leftover duplicate "Contact" row entirely, which shifts
# every following row (Description, Purpose, Copyright, FHIR Version, Kind,
# Type, Base Definition, Abstract, Derivation, Context) up by one.
$meta.Rows.Item(11).Delete()

$elements = $wb.Worksheets.Item(2)

# The root "Extension" element row: give it the real short/definition text.
$elements.Range("K2").Value = "RX Count"
$elements.Range("L2").Value = "Count of prescriptions for the drug claim"
